# 597_result.docx edit: split several <w:br/> line breaks into real
# paragraph breaks (new Heading-less paragraphs), and fix a handful of
# typos / word substitutions scattered through the same text.
#
# Strategy: the manual line breaks show up in Range.Text as char code 11
# (vertical tab). We locate each one with Find (narrow, unique literal
# context on both sides) and replace it with a real paragraph mark
# ("^p"), which splits the run/paragraph the same way Word does when you
# select the break and press Enter instead. Plain word-level fixes are
# done with ordinary literal Find & Replace.

$d = $word.ActiveDocument
$vt = [char]11   # the character a <w:br/> manual line break renders as

function Split-AtBreak($searchText) {
    # $searchText must contain exactly one $vt and have enough literal
    # context around it to be unique in the document.
    $replaceText = $searchText.Replace($vt, "^p")
    $ok = $d.Content.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Split-AtBreak: not found -> $searchText"
    }
}

function Replace-Text($findText, $replaceText) {
    $ok = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Replace-Text: not found -> $findText"
    }
}

# --- 1) Turn specific manual line breaks into paragraph breaks ---------

# "...(Asteraceae)." / "Tên khác : Pyrethre - Chrysanthème..."
Split-AtBreak ("Asteraceae)." + $vt + "Tên khác")

# "...Pyrethrum carneum M.B." / "Mô tả cây : Là cỏ sống dai..."
Split-AtBreak ("carneum M.B." + $vt + "Mô tả cây")

# "...chiếm 3/4 sản lượng." / "Thu hái chế biến : Gieo hạt..."
Split-AtBreak ("sản lượng." + $vt + "Thu hái chế biến")

# "...ngắt hoa đem chế biến." / "Thành phần hóa học : Có những chất :"
Split-AtBreak ("đem chế biến." + $vt + "Thành phần hóa học")

# "...tỷ lệ lại ít hơn (2/3)" / "Công dụng : Dùng để trừ sâu rau..."
Split-AtBreak ("ít hơn (2/3)" + $vt + "Công dụng")

# "...50 phần bột - nhựa làm hương." / "Lưu ý : Cúc trừ sâu gây ngộ độc..."
Split-AtBreak ("làm hương." + $vt + "Lưu ý")

# "...họ Đậu (Fabaceae)" / "Tên khác : Sắn nước..."
Split-AtBreak ("(Fabaceae)" + $vt + "Tên khác")

# "...Krásang (Campuchia)" / "Bộ phận dùng : Củ tươi..."
Split-AtBreak ("(Campuchia)" + $vt + "Bộ phận dùng")

# "...chưa dùng làm thuốc." / "Mô tả cây : Cây củ đậu..."
Split-AtBreak ("làm thuốc." + $vt + "Mô tả cây")

# --- 2) Word-level text fixes -------------------------------------------

Replace-Text "Chrysanthème vermicide" "Chrysanthène vermicide"
Replace-Text "phiến lá sẻ lông chim cắt sâu" "phiến lá xẻ lông chim cắt sâu"
Replace-Text "bên cạnh sẻ thùy sâu" "bên cạnh xẻ thùy sâu"
Replace-Text "hoa màu bị sâu phá hoại" "hoa mầu bị sâu phá hoại"
Replace-Text "bột thân lá cúc" "bột than lá cúc"
Replace-Text "bột cúc trừ sâu - bảo quản" "bột cúc trừ sâu : bảo quản"
Replace-Text "Củ sắn - Mănphầu (Lào)" "Củ sắng - Mănphău (Lào)"

Write-Host "Done. Paragraphs count:" $d.Paragraphs.Count
